# Update the F (and G) column numeric values on the "展览" and "全部类型"
# sheets, as produced by the regenerated data output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 819
$ws1.Range("F4").Value  = 1174
$ws1.Range("F5").Value  = 17
$ws1.Range("F11").Value = 528
$ws1.Range("F12").Value = 555
$ws1.Range("F13").Value = 165
$ws1.Range("F14").Value = 13174
$ws1.Range("G14").Value = 60
$ws1.Range("F18").Value = 5408
$ws1.Range("F19").Value = 5558
$ws1.Range("F20").Value = 21

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 819
$ws4.Range("F4").Value  = 1174
$ws4.Range("F12").Value = 17
$ws4.Range("F33").Value = 528
$ws4.Range("F34").Value = 555
$ws4.Range("F35").Value = 165
$ws4.Range("F36").Value = 13174
$ws4.Range("G36").Value = 60
$ws4.Range("F41").Value = 5408
$ws4.Range("F42").Value = 5558
$ws4.Range("F43").Value = 21
